$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="64.132.72"'
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Formula = '="3.482.71"'
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Formula = '="585.85"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Formula = '="132.16"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').Formula = '="7.64"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +4.77%  '
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Formula = '="4.075.71"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Formula = '="3.484.25"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Formula = '="64.106.38"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Formula = '="24.31"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -7.18%  '
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Formula = '="5.73"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Formula = '="13.55"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Formula = '="384.80"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').Formula = '="0.576"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Formula = '="3.623.32"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').Formula = '="74.74"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Formula = '="7.20"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -4.65%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Formula = '="1.44"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -3.78%  '
$ws.Range('D32').Formula = '="7.94"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  -4.09%  '
$ws.Range('E33').Value = '  +2.83%  '
$ws.Range('D34').Formula = '="3.512.96"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Formula = '="22.98"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('D37').Formula = '="5.19"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Formula = '="6.80"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('E39').Value = '  -3.29%  '
$ws.Range('D40').Formula = '="163.24"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +0.79%  '
$ws.Range('D41').Formula = '="0.0778"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('E44').Value = '  -1.60%  '
$ws.Range('D45').Formula = '="23.93"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -6.02%  '
$ws.Range('D46').Formula = '="1.63"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.32%  '
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('D48').Formula = '="0.923"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +3.18%  '
$ws.Range('D49').Formula = '="6.72"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').Formula = '="2.367.02"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('E51').Value = '  -2.41%  '
